$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6232884480899656
$ws.Range("D2").Value = 0.02690862442339181
$ws.Range("E2").Value = 0.2084112958037858
$ws.Range("F2").Value = 0.6917086417030944
$ws.Range("G2").Value = 0.002422417213734766
$ws.Range("I2").Value = 0.8987914908940198
$ws.Range("K2").Value = 0.4171906485533441
$ws.Range("L2").Value = 0.1735204502681711
$ws.Range("N2").Value = 1.606075040866205
$ws.Range("O2").Value = 2.359832708075601

$ws.Range("B3").Value = 0.5977414858373606
$ws.Range("D3").Value = 0.0250182463317401
$ws.Range("E3").Value = 0.2099774247134594
$ws.Range("F3").Value = 0.6879759794752189
$ws.Range("G3").Value = 0.002424782780578028
$ws.Range("I3").Value = 0.907204706526052
$ws.Range("K3").Value = 0.3639223348889971
$ws.Range("L3").Value = 0.1623285125416345
$ws.Range("N3").Value = 1.622308216799819
$ws.Range("O3").Value = 2.360461345777935

$ws.Range("B4").Value = 0.5822958695243585
$ws.Range("D4").Value = 0.02384711740899803
$ws.Range("E4").Value = 0.2110098799531901
$ws.Range("F4").Value = 0.6861100133191371
$ws.Range("G4").Value = 0.002426314196749564
$ws.Range("I4").Value = 0.9127680805732936
$ws.Range("K4").Value = 0.3310970610297659
$ws.Range("L4").Value = 0.1555327897606134
$ws.Range("N4").Value = 1.63280354603296
$ws.Range("O4").Value = 2.362229328441188

$ws.Range("B5").Value = 0.5760626450020254
$ws.Range("D5").Value = 0.02336727240821546
$ws.Range("E5").Value = 0.2114484587196541
$ws.Range("F5").Value = 0.6854567536984462
$ws.Range("G5").Value = 0.002426958171625151
$ws.Range("I5").Value = 0.9151352268179451
$ws.Range("K5").Value = 0.317691425004142
$ws.Range("L5").Value = 0.1527827495160494
$ws.Range("N5").Value = 1.637213304494603
$ws.Range("O5").Value = 2.363297343136821

$ws.Range("B6").Value = 0.5750313228724622
$ws.Range("D6").Value = 0.02328743809773215
$ws.Range("E6").Value = 0.2115223630618388
$ws.Range("F6").Value = 0.6853547524073917
$ws.Range("G6").Value = 0.002427066307281703
$ws.Range("I6").Value = 0.9155343321595133
$ws.Range("K6").Value = 0.3154636936105817
$ws.Range("L6").Value = 0.152327274874736
$ws.Range("N6").Value = 1.637953564941121
$ws.Range("O6").Value = 2.36349567817571

$ws.Range("B7").Value = 0.582211558393027
$ws.Range("D7").Value = 0.02384065655215295
$ws.Range("E7").Value = 0.2110157224860485
$ws.Range("F7").Value = 0.6861007693866412
$ws.Range("G7").Value = 0.002426322800886691
$ws.Range("I7").Value = 0.912799599680163
$ws.Range("K7").Value = 0.3309163847375771
$ws.Range("L7").Value = 0.1554956235684131
$ws.Range("N7").Value = 1.63286247972332
$ws.Range("O7").Value = 2.362242324813224

$ws.Range("B8").Value = 0.6144303027969045
$ws.Range("D8").Value = 0.02625899933916287
$ws.Range("E8").Value = 0.2089366111650692
$ws.Range("F8").Value = 0.6903332442784489
$ws.Range("G8").Value = 0.002423216510448967
$ws.Range("I8").Value = 0.9016098648018911
$ws.Range("K8").Value = 0.3988486828486373
$ws.Range("L8").Value = 0.1696457276739807
$ws.Range("N8").Value = 1.611562655342009
$ws.Range("O8").Value = 2.359762642155488

$ws.Range("B9").Value = 0.6794970697685017
$ws.Range("D9").Value = 0.03091784391224905
$ws.Range("E9").Value = 0.2054202593495162
$ws.Range("F9").Value = 0.7020121814853866
$ws.Range("G9").Value = 0.00241774885139025
$ws.Range("I9").Value = 0.8828199843867495
$ws.Range("K9").Value = 0.5311011740609501
$ws.Range("L9").Value = 0.1979947713493004
$ws.Range("N9").Value = 1.573981896577173
$ws.Range("O9").Value = 2.365867329335003

$ws.Range("B10").Value = 0.7284284918052322
$ws.Range("D10").Value = 0.0342890451585447
$ws.Range("E10").Value = 0.2031768073700828
$ws.Range("F10").Value = 0.7126544477688697
$ws.Range("G10").Value = 0.002414108310903015
$ws.Range("I10").Value = 0.8709347515764705
$ws.Range("K10").Value = 0.6276573783410697
$ws.Range("L10").Value = 0.2191863574853983
$ws.Range("N10").Value = 1.548920868109274
$ws.Range("O10").Value = 2.377045713747663

$ws.Range("B11").Value = 0.7509288231973699
$ws.Range("D11").Value = 0.03581133577122841
$ws.Range("E11").Value = 0.2022296404301738
$ws.Range("F11").Value = 0.7179440894775269
$ws.Range("G11").Value = 0.002412533108740953
$ws.Range("I11").Value = 0.8659441622722497
$ws.Range("K11").Value = 0.6714465664937848
$ws.Range("L11").Value = 0.2289054558372641
$ws.Range("N11").Value = 1.538072803982798
$ws.Range("O11").Value = 2.383586301068135

$ws.Range("B12").Value = 0.7594832957352935
$ws.Range("D12").Value = 0.03638614565642939
$ws.Range("E12").Value = 0.2018814961462621
$ws.Range("F12").Value = 0.7200116174371942
$ws.Range("G12").Value = 0.002411948194311338
$ws.Range("I12").Value = 0.8641141510310391
$ws.Range("K12").Value = 0.6880084420313324
$ws.Range("L12").Value = 0.2325970862033131
$ws.Range("N12").Value = 1.534044322518001
$ws.Range("O12").Value = 2.386272409666134

$ws.Range("B13").Value = 0.7576394305579015
$ws.Range("D13").Value = 0.03626242379082356
$ws.Range("E13").Value = 0.2019560074598044
$ws.Range("F13").Value = 0.7195634720685575
$ws.Range("G13").Value = 0.002412073651963498
$ws.Range("I13").Value = 0.864505616535947
$ws.Range("K13").Value = 0.6844424537583791
$ws.Range("L13").Value = 0.2318015307165666
$ws.Range("N13").Value = 1.5349083943729
$ws.Range("O13").Value = 2.385684597599493

$ws.Range("B14").Value = 0.7516319250969161
$ws.Range("D14").Value = 0.03585865891267304
$ws.Range("E14").Value = 0.2022007875577643
$ws.Range("F14").Value = 0.7181128946232604
$ws.Range("G14").Value = 0.002412484755961018
$ws.Range("I14").Value = 0.8657924074042569
$ws.Range("K14").Value = 0.6728095299929748
$ws.Range("L14").Value = 0.2292089444113969
$ws.Range("N14").Value = 1.537739786059813
$ws.Range("O14").Value = 2.38380309336921

$ws.Range("B15").Value = 0.7479565757564899
$ws.Range("D15").Value = 0.03561112581587622
$ws.Range("E15").Value = 0.2023520926222897
$ws.Range("F15").Value = 0.7172327672315362
$ws.Range("G15").Value = 0.002412738075241002
$ws.Range("I15").Value = 0.8665883934809493
$ws.Range("K15").Value = 0.6656813827063388
$ws.Range("L15").Value = 0.2276223675017803
$ws.Range("N15").Value = 1.539484441988066
$ws.Range("O15").Value = 2.382677879611691

$ws.Range("B16").Value = 0.7269628444638272
$ws.Range("D16").Value = 0.03418933059702312
$ws.Range("E16").Value = 0.203240181522192
$ws.Range("F16").Value = 0.712317778348563
$ws.Range("G16").Value = 0.002414212877412285
$ws.Range("I16").Value = 0.8712692694821911
$ws.Range("K16").Value = 0.6247928754961549
$ws.Range("L16").Value = 0.2185527686951474
$ws.Range("N16").Value = 1.549640930555414
$ws.Range("O16").Value = 2.376647561942917

$ws.Range("B17").Value = 0.7141452186077117
$ws.Range("D17").Value = 0.03331419650363188
$ws.Range("E17").Value = 0.203803773722262
$ws.Range("F17").Value = 0.7094174245307627
$ws.Range("G17").Value = 0.002415138303064886
$ws.Range("I17").Value = 0.874247389869705
$ws.Range("K17").Value = 0.599674051118285
$ws.Range("L17").Value = 0.2130089877765329
$ws.Range("N17").Value = 1.556013085124075
$ws.Range("O17").Value = 2.373320966184707

$ws.Range("B18").Value = 0.7067955960931442
$ws.Range("D18").Value = 0.03280978179008542
$ws.Range("E18").Value = 0.2041348464836847
$ws.Range("F18").Value = 0.7077914269630767
$ws.Range("G18").Value = 0.002415678200855535
$ws.Range("I18").Value = 0.8759994884909155
$ws.Range("K18").Value = 0.5852137337918748
$ws.Range("L18").Value = 0.2098277897659386
$ws.Range("N18").Value = 1.559730165733979
$ws.Range("O18").Value = 2.371544591492437

$ws.Range("B19").Value = 0.7043110634384675
$ws.Range("D19").Value = 0.03263881423452375
$ws.Range("E19").Value = 0.2042481295336245
$ws.Range("F19").Value = 0.7072481424054899
$ws.Range("G19").Value = 0.002415862310850361
$ws.Range("I19").Value = 0.876599445572726
$ws.Range("K19").Value = 0.5803155723109796
$ws.Range("L19").Value = 0.2087519723166196
$ws.Range("N19").Value = 1.560997633365133
$ws.Range("O19").Value = 2.370966669555997

$ws.Range("B20").Value = 0.7155073279007524
$ws.Range("D20").Value = 0.03340746598843936
$ws.Range("E20").Value = 0.2037430634604593
$ws.Range("F20").Value = 0.7097218038395638
$ws.Range("G20").Value = 0.002415039001773145
$ws.Range("I20").Value = 0.8739263107863273
$ws.Range("K20").Value = 0.6023493059781799
$ws.Range("L20").Value = 0.2135983638580541
$ws.Range("N20").Value = 1.555329377903591
$ws.Range("O20").Value = 2.37366090924229

$ws.Range("B21").Value = 0.7533955547624771
$ws.Range("D21").Value = 0.03597729939033201
$ws.Range("E21").Value = 0.2021286042103867
$ws.Range("F21").Value = 0.7185372153007847
$ws.Range("G21").Value = 0.002412363690814759
$ws.Range("I21").Value = 0.8654128224287447
$ws.Range("K21").Value = 0.6762269543982597
$ws.Range("L21").Value = 0.2299701464422697
$ws.Range("N21").Value = 1.536905981388793
$ws.Range("O21").Value = 2.384350055356492

$ws.Range("B22").Value = 0.7783560998717576
$ws.Range("D22").Value = 0.03764721320710152
$ws.Range("E22").Value = 0.2011348107117534
$ws.Range("F22").Value = 0.7246742785116282
$ws.Range("G22").Value = 0.002410682694308885
$ws.Range("I22").Value = 0.8601974168527597
$ws.Range("K22").Value = 0.7243925220126357
$ws.Range("L22").Value = 0.2407353820970144
$ws.Range("N22").Value = 1.525328304723043
$ws.Range("O22").Value = 2.392556180275363

$ws.Range("B23").Value = 0.7650162308131598
$ws.Range("D23").Value = 0.03675683758354609
$ws.Range("E23").Value = 0.2016596122049634
$ws.Range("F23").Value = 0.7213644466831539
$ws.Range("G23").Value = 0.002411573717012475
$ws.Range("I23").Value = 0.8629490795192325
$ws.Range("K23").Value = 0.6986966798646677
$ws.Range("L23").Value = 0.2349838396378487
$ws.Range("N23").Value = 1.531465149208294
$ws.Range("O23").Value = 2.38806476275775

$ws.Range("B24").Value = 0.7148914578118024
$ws.Range("D24").Value = 0.03336530286600947
$ws.Range("E24").Value = 0.2037704885895817
$ws.Range("F24").Value = 0.7095840648244405
$ws.Range("G24").Value = 0.002415083871422218
$ws.Range("I24").Value = 0.8740713462364909
$ws.Range("K24").Value = 0.6011398829133157
$ws.Range("L24").Value = 0.2133318882470974
$ws.Range("N24").Value = 1.555638314840401
$ws.Range("O24").Value = 2.373506796952398

$ws.Range("B25").Value = 0.6616952457260936
$ws.Range("D25").Value = 0.02966652012170101
$ws.Range("E25").Value = 0.2063116824795053
$ws.Range("F25").Value = 0.6984907791853772
$ws.Range("G25").Value = 0.002419161605275004
$ws.Range("I25").Value = 0.8875658429936948
$ws.Range("K25").Value = 0.4954287837333311
$ws.Range("L25").Value = 0.2318015307165666
$ws.Range("N25").Value = 1.583700506147347
$ws.Range("O25").Value = 2.363040961651876
